$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.928.11"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.183.33"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.69"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.25"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.403"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.736"
$ws.Range("E8").Value = "  +3.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.185.08"
$ws.Range("E10").Value = "  -0.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.569"
$ws.Range("E11").Value = "  +1.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.183"
$ws.Range("E12").Value = "  +2.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.605.99"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.34"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.760.96"
$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.63"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.182.40"
$ws.Range("E18").Value = "  -0.88%  "

$ws.Range("E19").Value = "  +3.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000214"
$ws.Range("E20").Value = "  +31.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "434.83"
$ws.Range("E22").Value = "  +3.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.49"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("E24").Value = "  -4.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.31"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.73"
$ws.Range("E26").Value = "  -5.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "80.98"
$ws.Range("E27").Value = "  +8.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.347.10"
$ws.Range("E28").Value = "  +0.40%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  -6.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.06"
$ws.Range("E32").Value = "  +27.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.38"
$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "514.88"
$ws.Range("E34").Value = "  -8.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.01"
$ws.Range("E35").Value = "  -1.01%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("E37").Value = "  -4.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.35"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.36"
$ws.Range("E39").Value = "  +2.37%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.127"
$ws.Range("E41").Value = "  -4.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.92"
$ws.Range("E43").Value = "  -1.77%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.373"
$ws.Range("E44").Value = "  -2.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.28"
$ws.Range("E45").Value = "  -2.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.96"
$ws.Range("E46").Value = "  +1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "169.36"
$ws.Range("E47").Value = "  -4.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.740"
$ws.Range("E49").Value = "  +5.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.81"
$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("E51").Value = "  -4.15%  "
